# Insert a new weekly data row at row 314 (pushing the previous rows 314-410 down to 315-411),
# then populate the newly inserted row 314 with the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 314. Excel shifts rows 314:410 down to 315:411.
$ws.Rows("314:314").Insert()

# Populate the new row 314 with the new data record.
$ws.Cells.Item(314, 1).Value = 6
$ws.Cells.Item(314, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(314, 3).Value = "Metropolitana"
$ws.Cells.Item(314, 4).Value = 44736
$ws.Cells.Item(314, 5).Value = 13
$ws.Cells.Item(314, 6).Value = 100112032
$ws.Cells.Item(314, 7).Value = "Zapallo italiano"
$ws.Cells.Item(314, 8).Value = "Sin especificar"
$ws.Cells.Item(314, 9).Value = "Primera"
$ws.Cells.Item(314, 10).Value = 580
$ws.Cells.Item(314, 11).Value = 13000
$ws.Cells.Item(314, 12).Value = 15000
$ws.Cells.Item(314, 13).Value = 13897
$ws.Cells.Item(314, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(314, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(314, 16).Value = 278
$ws.Cells.Item(314, 17).Value = 50
$ws.Cells.Item(314, 18).Value = "Hortaliza"
